# Remove the first paragraph of the document body (the italic
# "I personally examined the patient..." attestation paragraph),
# including its paragraph mark, so the following "OBJECTIVE:"
# paragraph becomes the first paragraph in the document.

$d = $word.ActiveDocument
$d.Paragraphs(1).Range.Delete()
